# Update countries & provincias Spain
# Applies the data refresh to the "Pais" worksheet:
#  - updates the "Datos actualizados..." timestamp string
#  - updates case-count figures for several countries
#  - Georgia overtakes Senegal (rows 111/112 swap country + figures)
#  - Uganda overtakes Islas Caimanes (rows 149/150 swap country + figures)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 23:52"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 610467
$ws.Range("C4").Value = 23526
$ws.Range("E4").Value = 546093
$ws.Range("G4").Value = 2214
$ws.Range("H4").Value = 25854

# --- Row 29: Ecuador ---
$ws.Range("E29").Value = 6538
$ws.Range("F29").Value = 129
$ws.Range("G29").Value = 14
$ws.Range("H29").Value = 369

# --- Rows 111/112: Georgia overtakes Senegal ---
$ws.Range("A111").Value = "Georgia"
$ws.Range("B111").Value = 300
$ws.Range("C111").Value = 28
$ws.Range("D111").Value = 69
$ws.Range("E111").Value = 228
$ws.Range("F111").Value = 6
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 3

$ws.Range("A112").Value = "Senegal"
$ws.Range("B112").Value = 299
$ws.Range("C112").Value = 8
$ws.Range("D112").Value = 183
$ws.Range("E112").Value = 114
$ws.Range("F112").Value = 1
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 2

# --- Rows 149/150: Uganda overtakes Islas Caimanes ---
$ws.Range("A149").Value = "Uganda"
$ws.Range("B149").Value = 55
$ws.Range("C149").Value = 1
$ws.Range("D149").Value = 8
$ws.Range("E149").Value = 47
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 0

$ws.Range("A150").Value = "Islas Caimanes"
$ws.Range("B150").Value = 54
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 6
$ws.Range("E150").Value = 47
$ws.Range("F150").Value = 3
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 1
